# y20d02 results: fill in "Day 2: Password Philosophy" row and re-point selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2020")

# Row 6 used to be a placeholder "Day 2: " row with no data; give it the
# real puzzle title and this year's timings/rank.
$ws.Range("B6").Value = "Day 2: Password Philosophy"
$ws.Range("C6").Value = 0.17916666666666667
$ws.Range("E6").Value = 0.29236111111111113
$ws.Range("F6").Value = 0.25694444444444448
$ws.Range("H6").Value = "2nd"

# D6:D29 mirror D5's "finish - start" elapsed-time formula (fill down as one
# relative-formula range so Excel keeps it as a single shared formula group).
$ws.Range("D6:D29").Formula = "=E6-C6"

# Move the selection to reflect where editing left off.
$ws.Range("H7").Select()
